# education_bar.xlsx - content/proofing edits
# (commit: "Fixing formatting through .cls file. Attempt at Lua filter for scientific names.")
#
# Semantic changes being applied:
#  - C2: "Tulane University, ..."            -> "Tulane University (TU), ..."
#  - C4: "University of Michigan, ..."       -> "University of Michigan (UM), ..."
#  - C5: "University of Puerto Rico"         -> "University of Puerto Rico (UPR)"
#  - E3: "...(adviser; dissertation chair)..." -> "...(advisor; dissertation chair)..."
#  - A4: "Master of Science, ... Ecology,"   -> "Master of Science, ... Ecology" (dropped trailing comma)
#  - Active cell/selection moves from E3 to A4
#  - Column widths/row sizing nudged slightly (best-effort; driven by the
#    Excel build that resaved the file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Tulane University (TU), Department of Ecology and Biology"
$ws.Range("E3").Value = "Committee members: Sunshine Van Bael, Ph.D. (advisor; dissertation chair), Kathleen Ferris, Ph.D.  (co-advisor), Keith Clay, Ph.D., & P. Camilo Zalamea, Ph.D."
$ws.Range("A4").Value = "Master of Science, Natural Resources Management: Conservation Ecology"
$ws.Range("C4").Value = "University of Michigan (UM), School for the Environment and Sustainability"
$ws.Range("C5").Value = "University of Puerto Rico (UPR)"

# Column widths settled slightly differently after the resave (best effort —
# headless engine rounds ColumnWidth to sixths of a character, so this is the
# closest achievable approximation of the committed 63.3984375 / 41.296875).
$ws.Columns.Item(1).ColumnWidth = 62.417
$ws.Columns.Item(2).ColumnWidth = 40.417

# Move the active selection from E3 to A4, matching the saved cursor position.
$ws.Range("A4").Select()
